$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 12:29"

# Update country rows: new countries inserted into the ranking shift other rows down,
# and daily case/death counts were refreshed for the new snapshot.
$ws.Range("E8").Value = 46283
$ws.Range("G8").Value = 48
$ws.Range("H8").Value = 399
$ws.Range("B20").Value = 3807
$ws.Range("C20").Value = 36
$ws.Range("E20").Value = 3780
$ws.Range("A33").Value = "Rumania"
$ws.Range("B33").Value = 1452
$ws.Range("C33").Value = 160
$ws.Range("D33").Value = 139
$ws.Range("E33").Value = 1284
$ws.Range("F33").Value = 34
$ws.Range("G33").Value = 3
$ws.Range("H33").Value = 29
$ws.Range("A34").Value = "Polonia"
$ws.Range("B34").Value = 1436
$ws.Range("C34").Value = 47
$ws.Range("D34").Value = 7
$ws.Range("E34").Value = 1413
$ws.Range("F34").Value = 3
$ws.Range("H34").Value = 16
$ws.Range("A35").Value = "Pakistan"
$ws.Range("B35").Value = 1408
$ws.Range("C35").Value = 35
$ws.Range("D35").Value = 25
$ws.Range("E35").Value = 1372
$ws.Range("F35").Value = 7
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 11
$ws.Range("A61").Value = "Irak"
$ws.Range("B61").Value = 506
$ws.Range("C61").Value = 48
$ws.Range("D61").Value = 131
$ws.Range("E61").Value = 333
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 42
$ws.Range("A62").Value = "Barein"
$ws.Range("B62").Value = 473
$ws.Range("C62").Value = 7
$ws.Range("D62").Value = 254
$ws.Range("E62").Value = 215
$ws.Range("F62").Value = 1
$ws.Range("H62").Value = 4
$ws.Range("B98").Value = 120
$ws.Range("C98").Value = 5
$ws.Range("D98").Value = 25
$ws.Range("E98").Value = 94
$ws.Range("B101").Value = 110
$ws.Range("C101").Value = 4
$ws.Range("E101").Value = 101
$ws.Range("A119").Value = "Consejo Danes para los Refugiados"
$ws.Range("B119").Value = 58
$ws.Range("C119").Value = 7
$ws.Range("D119").Value = 2
$ws.Range("E119").Value = 50
$ws.Range("G119").Value = 3
$ws.Range("H119").Value = 6
$ws.Range("A120").Value = "Liechtenstein"
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 0
$ws.Range("E120").Value = 56
$ws.Range("F120").Value = 0
$ws.Range("H120").Value = 0
$ws.Range("A121").Value = "Paraguay"
$ws.Range("B121").Value = 56
$ws.Range("C121").Value = 4
$ws.Range("D121").Value = 1
$ws.Range("E121").Value = 52
$ws.Range("F121").Value = 1
$ws.Range("H121").Value = 3
$ws.Range("A122").Value = "Gibraltar"
$ws.Range("B122").Value = 55
$ws.Range("D122").Value = 14
$ws.Range("E122").Value = 41
$ws.Range("A123").Value = "Ruanda"
$ws.Range("B123").Value = 54
$ws.Range("D123").Value = 0
$ws.Range("E123").Value = 54
$ws.Range("H123").Value = 0
$ws.Range("A152").Value = "San Martin (Parte Francesa)"
$ws.Range("A153").Value = "Mali"
$ws.Range("A154").Value = "Dominica"
$ws.Range("A159").Value = "Surinam"
$ws.Range("A161").Value = "Haiti"
$ws.Range("A166").Value = "Seychelles"
$ws.Range("A167").Value = "Antigua y Barbuda"
$ws.Range("A168").Value = "Granada"
$ws.Range("A169").Value = "Mozambique"
$ws.Range("A170").Value = "Zimbabue"
$ws.Range("C170").Value = 2
$ws.Range("A171").Value = "Gabon"
$ws.Range("B171").Value = 7
$ws.Range("H171").Value = 1
$ws.Range("A172").Value = "Laos"
$ws.Range("A174").Value = "Benin"
$ws.Range("B174").Value = 6
$ws.Range("E174").Value = 6
$ws.Range("A175").Value = "Mauritania"
$ws.Range("C175").Value = 2
$ws.Range("A176").Value = "San Bartolome"
$ws.Range("C176").Value = 0
$ws.Range("A177").Value = "Fiyi"
$ws.Range("A178").Value = "Siria"
$ws.Range("A179").Value = "Montserrat"
$ws.Range("E179").Value = 5
$ws.Range("H179").Value = 0
$ws.Range("A181").Value = "Cabo Verde"
$ws.Range("A182").Value = "Guyana"
$ws.Range("C182").Value = 0
$ws.Range("A183").Value = "Sudan"
$ws.Range("C183").Value = 2
$ws.Range("A184").Value = "Congo"
$ws.Range("A185").Value = "Angola"
$ws.Range("A188").Value = "Butan"
$ws.Range("A189").Value = "Liberia"
$ws.Range("A190").Value = "Somalia"
$ws.Range("A191").Value = "Republica del Chad"
